$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("commondata")

# Clear the password value in B3 (the "123456" shared string) while
# keeping the cell's existing style (quote-prefix style stays intact).
$ws.Range("B3").ClearContents()
